# Applies the "Updated cryptos list" refresh described by the commit diff.
# Row order: Solana/Cardano (rows 8-9) and FraxShare/PaxDollar (rows 41-42) swap places,
# prices and 1h volume percentages are refreshed across the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.043.06'
$ws.Range("E2").Value = '  +0.64%  '
$ws.Range("D3").Value = '1.680.72'
$ws.Range("E3").Value = '  +0.87%  '
$ws.Range("D5").Value = "'215.96"
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("E6").Value = '  -3.08%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = "'0.254"
$ws.Range("E8").Value = '  +1.20%  '
$ws.Range("B9").Value = 'Solana'
$ws.Range("C9").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D9").Value = "'21.56"
$ws.Range("E9").Value = '  +6.83%  '
$ws.Range("E10").Value = '  +0.59%  '
$ws.Range("D11").Value = "'0.0889"
$ws.Range("E11").Value = '  -0.74%  '
$ws.Range("D12").Value = '1.917.08'
$ws.Range("E12").Value = '  +0.81%  '
$ws.Range("D13").Value = '1.673.59'
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("D14").Value = "'4.11"
$ws.Range("E14").Value = '  +0.56%  '
$ws.Range("E15").Value = '  +1.57%  '
$ws.Range("D16").Value = "'66.53"
$ws.Range("E16").Value = '  +0.76%  '
$ws.Range("D17").Value = '27.033.25'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("D18").Value = "'8.18"
$ws.Range("E18").Value = '  +4.83%  '
$ws.Range("D19").Value = "'235.91"
$ws.Range("E19").Value = '  +1.64%  '
$ws.Range("D20").Value = '0.0₃0739'
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("E23").Value = '  +1.05%  '
$ws.Range("E24").Value = '  -4.18%  '
$ws.Range("D25").Value = "'146.69"
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("D26").Value = "'7.25"
$ws.Range("E26").Value = '  +1.82%  '
$ws.Range("D27").Value = "'16.45"
$ws.Range("E27").Value = '  +3.60%  '
$ws.Range("E28").Value = '  -2.62%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  +0.62%  '
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("E32").Value = '  -0.15%  '
$ws.Range("D33").Value = '1.525.73'
$ws.Range("E33").Value = '  +4.21%  '
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("E35").Value = '  +5.04%  '
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("E37").Value = '  +3.43%  '
$ws.Range("D38").Value = "'0.922"
$ws.Range("E38").Value = '  +2.57%  '
$ws.Range("E39").Value = '  +3.10%  '
$ws.Range("E40").Value = '  +5.45%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = "'5.68"
$ws.Range("E42").Value = '  -3.02%  '
$ws.Range("D43").Value = "'67.93"
$ws.Range("E43").Value = '  +3.13%  '
$ws.Range("E44").Value = '  -0.34%  '
$ws.Range("D45").Value = '1.822.07'
$ws.Range("E45").Value = '  +0.23%  '
$ws.Range("E46").Value = '  +0.40%  '
$ws.Range("D47").Value = "'90.31"
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("E48").Value = '  -0.14%  '
$ws.Range("E49").Value = '  +2.48%  '
$ws.Range("D50").Value = "'7.92"
$ws.Range("E50").Value = '  +4.32%  '
$ws.Range("E51").Value = '  -0.47%  '
